{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The first paragraph is \"Test Document for Document Processor\".\nconst firstParagraph = paragraphs.items[0];\n\n// Insert the three new paragraphs right after it, in order.\nconst emptyParagraph = firstParagraph.insertParagraph(\"\", \"After\");\n// Force an explicit (empty) <w:t> run of text so the paragraph matches the\n// source document's convention for blank paragraphs.\nemptyParagraph.getRange().insertText(\"\", \"Replace\");\nconst headingParagraph = emptyParagraph.insertParagraph(\"Grammar Test Paragraph:\", \"After\");\nheadingParagraph.insertParagraph(\n  \"The cats and dog is running fast. We dont need no help with grammer. This sentense contains muliple mispelled words. The weather have been nice yesterday?\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The first paragraph is \"Test Document for Document Processor\".\n$p1 = $d.Paragraphs(1)\n\n# Insert a new (blank) paragraph right after it.\n$p1.Range.InsertParagraphAfter()\n\n# That new blank paragraph is now paragraph 2. Force an explicit (empty)\n# text run so it matches the source document's convention for blank\n# paragraphs (an empty <w:t/> rather than a completely empty run).\n$p2 = $d.Paragraphs(2)\n$p2.Range.Text = \"\"\n\n# Insert another new paragraph after it, and give it the heading text.\n$p2.Range.InsertParagraphAfter()\n$p3 = $d.Paragraphs(3)\n$p3.Range.InsertBefore(\"Grammar Test Paragraph:\")\n\n# Insert the final new paragraph (the grammar-error sentence) after that.\n$p3.Range.InsertParagraphAfter()\n$p4 = $d.Paragraphs(4)\n$p4.Range.InsertBefore(\"The cats and dog is running fast. We dont need no help with grammer. This sentense contains muliple mispelled words. The weather have been nice yesterday?\")\n"}
